$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Map of shape index -> new text for the second paragraph (drop the trailing
# numeric suffix that used to follow the closing bracket, e.g. "[RbtWeld01]2"
# -> "[RbtWeld01]").
$changes = @{
    3  = "[RbtWeld01]"
    4  = "[2ndClamp01]"
    5  = "[1stClamp01]"
    14 = "[RbtWeld02]"
    15 = "[2ndClamp02]"
}

foreach ($idx in $changes.Keys) {
    $sh = $s.Shapes.Item($idx)
    $tr = $sh.TextFrame.TextRange
    $para2 = $tr.Paragraphs(2,1)
    $para2.Text = $changes[$idx]
}
